# chilkat_sockeye_template.xlsx update
# - Update review comments (DIDSON data range extended through 2022,
#   clarified expansion-factor guidance, added note about whole numbers
#   in the brood-table "age" sheet / Total Run table).
# - Refresh the age-composition percentages (columns J/K/L) for the most
#   recent run years now that 2022 data is final, and clear the leftover
#   "example output" highlighting / flagged formatting on the 2017-2022 rows.
# - Move the active selection to L22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Comment text updates
# ---------------------------------------------------------------------

$commentB1 = $ws.Range("B1").Comment
$commentB1.Text("SEM:; 3-17-2023; Make sure the spreadsheet Chilkat Lake DIDSON-daily counts 2008-2022 is updated to reflect any DIDSON expansion. This will carry over to the 'Chikat sockeye weighted age comp.' spreadsheet and the brood table spreadsheet, 'age' sheet. `n")

$commentC1 = $ws.Range("C1").Comment
$commentC1.Text("SEM: 3/20/2023 verify 0.05 if not expanded; 0.1 if expanded early OR late; 0.2 if expanded for both early and late`n")

$commentJ1 = $ws.Range("J1").Comment
$commentJ1.Text("SEM; 3-17-2023: From the brood table spreadsheet, 'age' tab. Look at table Total Run (effective sample size); make sure these are whole numbers!`n")

# ---------------------------------------------------------------------
# 2. Age composition value updates (columns J/K/L, various rows)
# ---------------------------------------------------------------------

$ws.Range("J2").Value = 6
$ws.Range("L2").Value = 28

$ws.Range("L3").Value = 21

$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 65
$ws.Range("L4").Value = 29

$ws.Range("K5").Value = 77
$ws.Range("L5").Value = 21

$ws.Range("L6").Value = 35

$ws.Range("L7").Value = 43

$ws.Range("K8").Value = 51
$ws.Range("L8").Value = 46

$ws.Range("K9").Value = 56
$ws.Range("L9").Value = 38

$ws.Range("K10").Value = 74
$ws.Range("L10").Value = 24

$ws.Range("K11").Value = 39
$ws.Range("L11").Value = 58

$ws.Range("J13").Value = 4
$ws.Range("K13").Value = 59

$ws.Range("K14").Value = 43
$ws.Range("L14").Value = 55

$ws.Range("K17").Value = 47
$ws.Range("L17").Value = 51

$ws.Range("L19").Value = 58

$ws.Range("L25").Value = 56

$ws.Range("L27").Value = 19

$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 38

$ws.Range("K32").Value = 60
$ws.Range("L32").Value = 34

$ws.Range("J33").Value = 6
$ws.Range("K33").Value = 52
$ws.Range("L33").Value = 42

$ws.Range("J34").Value = 6

$ws.Range("K41").Value = 29
$ws.Range("L41").Value = 67

$ws.Range("K42").Value = 59

$ws.Range("L48").Value = 12

# ---------------------------------------------------------------------
# 3. Clear the stray formatting (fill highlight / flagged-but-blank
#    fill) left over on the 2017-2022 rows from the example output.
# ---------------------------------------------------------------------

$ws.Range("B43:C48").ClearFormats()
$ws.Range("J43:L48").ClearFormats()

# ---------------------------------------------------------------------
# 4. Move the selection to L22
# ---------------------------------------------------------------------

$ws.Range("L22").Select()
